# Natmi following Dr Hou advice
# Rewrites the LR-pairs results grid (Sheet1) with the 3x3 cartesian
# product of sending/target clusters (ECs/FAPs/sCs) for Col1a1-Cd36,
# replacing the previous 6-row (incomplete) table with the full 9-row
# table and updated statistics.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns: A Sending cluster, B Ligand symbol, C Receptor symbol,
#          D Target cluster, E..T the numeric statistics.
$data = @(
    @("ECs","Col1a1","Cd36","ECs",3.0,1.0,6.072131,18.216393,0.003943999267036455,0.003943999267036454,3.0,1.0,550.056952,1650.170856,0.923678797551236,0.923678797551236,3340.017870004712,30060.16083004241,0.003642988500519189,0.003642988500519188),
    @("ECs","Col1a1","Cd36","FAPs",3.0,1.0,6.072131,18.216393,0.003943999267036455,0.003943999267036454,2.0,0.6666666666666666,1.448019333333334,4.344058,0.002431575042877152,0.002431575042877152,8.792563082532668,79.13306774279401,0.000009590130186851625,0.000009590130186851623),
    @("ECs","Col1a1","Cd36","sCs",3.0,1.0,6.072131,18.216393,0.003943999267036455,0.003943999267036454,3.0,1.0,44.00177133333333,132.005314,0.07388962740588682,0.07388962740588684,267.1845197680447,2404.660677912402,0.0002914206363304144,0.0002914206363304144),
    @("FAPs","Col1a1","Cd36","ECs",3.0,1.0,1480.851806666667,4442.55542,0.9618498744646554,0.9618498744646552,3.0,1.0,550.056952,1650.170856,0.923678797551236,0.923678797551236,814552.83113876,7330975.48024884,0.8884403354703203,0.8884403354703201),
    @("FAPs","Col1a1","Cd36","FAPs",3.0,1.0,1480.851806666667,4442.55542,0.9618498744646554,0.9618498744646552,2.0,0.6666666666666666,1.448019333333334,4.344058,0.002431575042877152,0.002431575042877152,2144.302045854929,19298.71841269436,0.002338810149742778,0.002338810149742778),
    @("FAPs","Col1a1","Cd36","sCs",3.0,1.0,1480.851806666667,4442.55542,0.9618498744646554,0.9618498744646552,3.0,1.0,44.00177133333333,132.005314,0.07388962740588682,0.07388962740588684,65160.1025755002,586440.9231795018,0.07107072884459241,0.0710707288445924),
    @("sCs","Col1a1","Cd36","ECs",3.0,1.0,52.663316,157.989948,0.03420612626830831,0.0342061262683083,3.0,1.0,550.056952,1650.170856,0.923678797551236,0.923678797551236,28967.82308117284,260710.4077305555,0.03159547358039677,0.03159547358039676),
    @("sCs","Col1a1","Cd36","FAPs",3.0,1.0,52.663316,157.989948,0.03420612626830831,0.0342061262683083,2.0,0.6666666666666666,1.448019333333334,4.344058,0.002431575042877152,0.002431575042877152,76.25749972544268,686.3174975289841,0.00008317476294752307,0.00008317476294752304),
    @("sCs","Col1a1","Cd36","sCs",3.0,1.0,52.663316,157.989948,0.03420612626830831,0.0342061262683083,3.0,1.0,44.00177133333333,132.005314,0.07388962740588682,0.07388962740588684,2317.279188287075,20855.51269458367,0.002527477924964019,0.002527477924964019)
)

$r = 2
foreach ($row in $data) {
    $c = 1
    foreach ($val in $row) {
        $ws.Cells.Item($r, $c).Value = $val
        $c = $c + 1
    }
    $r = $r + 1
}
